$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 499

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -ne $null) {
        $s = [string]$v
        if ($s.EndsWith("16")) {
            $cell.Value = $s.Substring(0, $s.Length - 2)
        }
    }
}
